$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-04-22 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-23 Tuesday", 2)

# Update the division problems table by direct cell addressing (row, col)
# since some values repeat with different replacements in different cells.
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cellRange = $table.Cell($row, $col).Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $text
}

Set-CellText $t 1 1 "44÷9="
Set-CellText $t 1 2 "88÷3="
Set-CellText $t 1 3 "69÷9="
Set-CellText $t 1 4 "67÷2="
Set-CellText $t 1 5 "19÷2="

Set-CellText $t 5 1 "15÷6="
Set-CellText $t 5 2 "60÷9="
Set-CellText $t 5 3 "56÷4="
Set-CellText $t 5 4 "85÷4="
Set-CellText $t 5 5 "26÷9="

Set-CellText $t 9 1 "12÷2="
Set-CellText $t 9 2 "29÷3="
Set-CellText $t 9 3 "56÷4="
Set-CellText $t 9 4 "41÷4="
Set-CellText $t 9 5 "33÷5="

Set-CellText $t 13 1 "89÷6="
Set-CellText $t 13 2 "64÷9="
Set-CellText $t 13 3 "31÷2="
Set-CellText $t 13 4 "41÷4="
Set-CellText $t 13 5 "41÷7="

Set-CellText $t 17 1 "23÷9="
Set-CellText $t 17 2 "63÷7="
Set-CellText $t 17 3 "58÷9="
Set-CellText $t 17 4 "67÷7="
Set-CellText $t 17 5 "38÷3="
